$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(50, 1).Formula = "9.199999999999999"
Write-Output ("val: " + $ws.Cells.Item(50,1).Value())
